# The workbook's "Perejil" (parsley) price log gains one new weekly data
# point. This is implemented as a row insert at row 163 (pushing the
# existing rows 163:266 down to 164:267, which is exactly what the
# target diff shows - every row's content equals what used to be one row
# above it), followed by filling in the brand-new row 163 with the new
# observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 163; Excel shifts 163:266 down to 164:267.
$ws.Rows("163:163").Insert()

# Populate the newly inserted row 163 with the new weekly record.
$ws.Range("A163").Value = 10
$ws.Range("B163").Value = 'Vega Modelo de Temuco'
$ws.Range("C163").Value = 'La Araucanía'
$ws.Range("D163").Value = 44596
$ws.Range("E163").Value = 9
$ws.Range("F163").Value = 100112044
$ws.Range("G163").Value = 'Perejil'
$ws.Range("H163").Value = 'Sin especificar'
$ws.Range("I163").Value = 'Primera'
$ws.Range("J163").Value = 40
$ws.Range("K163").Value = 5000
$ws.Range("L163").Value = 5000
$ws.Range("M163").Value = 5000
$ws.Range("N163").Value = '$/docena de atados (3 kilos)'
$ws.Range("O163").Value = 'Provincia de Cautín'
$ws.Range("P163").Value = 1667
$ws.Range("Q163").Value = 3
$ws.Range("R163").Value = 'Hortaliza'
